$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New log entries for 21-24 June 2022 (rows 19-23), plus the relocated
# "Total Hours" summary row (old row 22 -> new row 25, row 24 left blank).
# ---------------------------------------------------------------------------

function Add-LogRow {
    param($Row, $Sno, $DateSerial, $StartTime, $EndTime, $Category, $Description, $Height)

    $ws.Range("A$Row").Value = $Sno
    $ws.Range("B$Row").Value = $DateSerial
    $ws.Range("B$Row").NumberFormat = "m/d/yy"

    $ws.Range("C$Row").Value = $StartTime
    $ws.Range("C$Row").NumberFormat = "h:mm AM/PM"

    $ws.Range("D$Row").Value = $EndTime
    $ws.Range("D$Row").NumberFormat = "h:mm AM/PM"

    $ws.Range("E$Row").Formula = "=D$Row-C$Row"
    $ws.Range("E$Row").NumberFormat = "h:mm"

    $ws.Range("F$Row").Value = $Category

    $ws.Range("G$Row").Value = $Description
    $ws.Range("G$Row").WrapText = $true

    if ($Height) {
        $ws.Rows.Item($Row).RowHeight = $Height
    }
}

Add-LogRow 19 18 44733 0.29166666666666669 0.33333333333333331 "Code" "1. CE with weights vs Dice loss combinations`r`n2. deeplabv3_r50 vs r50 vs pt implementation`r`n3. different instances comparison" 45

Add-LogRow 20 19 44733 0.89583333333333337 0.92708333333333337 "Code" "1. Corrected FCN_r50 model with extra arguments" $null

Add-LogRow 21 20 44734 0.22916666666666666 0.28125 "Code" "1. common labels and images in BDD100k dataset`r`n2. resized images and upload in drive" 30

Add-LogRow 22 21 44735 0.16666666666666666 0.27083333333333331 "Code" "1. tried converting png to jpg, for some reason jpg and png are extremely slow to train`r`n2. resized images, labels to 320 x 180 and stored in npy format`r`n3. end to end training of FCN_resnet50_starter notebook" 75

Add-LogRow 23 22 44736 0.3125 0.36458333333333331 "Code" "1. created test video from BDD100k dataset`r`n2. Produced model output from input video" 30

# ---------------------------------------------------------------------------
# Relocate the "Total Hours" row from 22 to 25 (row 24 stays blank), and
# extend the SUM range along with it. (The old C22/E22 cells were already
# overwritten above by the new 23-June log entry, so nothing to clear here.)
# ---------------------------------------------------------------------------

$ws.Range("C25").Value = "Total Hours"

$ws.Range("E25").Formula = "=SUM(E2:E21)"
$ws.Range("E25").NumberFormat = "[hh]:mm"

# ---------------------------------------------------------------------------
# Window / selection state to mirror the edited view.
# ---------------------------------------------------------------------------

$ws.Range("G26").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$win.ScrollColumn = 1

$wb.Application.Calculate()
